$d = $word.ActiveDocument

# "Matriz Cuadrada de Transacciones" is the first table in the document.
$t = $d.Tables(1)

# Cells to mark with "X" (Word table is 1-indexed: Cell(row, column)).
# Row labels (entities) / Column = C-"L" or D-"L" sub-column of the
# square transactions matrix (A..E groups, each with I/L/A/B columns).
$cellsToFill = @(
    @{ Row = 5;  Col = 11 },  # ARTICULO  - C / L
    @{ Row = 6;  Col = 11 },  # CLIENTE   - C / L
    @{ Row = 14; Col = 11 },  # FACTURA   - C / L
    @{ Row = 17; Col = 11 },  # INCLUYE   - C / L
    @{ Row = 20; Col = 15 },  # ORDEN     - D / L  (also gets the _GoBack bookmark)
    @{ Row = 28; Col = 15 }   # USA       - D / L
)

# Move the stray "_GoBack" bookmark (currently sitting in its own empty
# paragraph right after the fourth table) onto the ORDEN/D-L cell so it
# ends up right after the new "X" run, matching a fresh Word edit there.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

foreach ($c in $cellsToFill) {
    $cell = $t.Cell($c.Row, $c.Col)

    if ($c.Row -eq 20 -and $c.Col -eq 15) {
        # Anchor the bookmark on the still-empty cell first so that the
        # text typed afterwards lands before it (matching the diff, where
        # the bookmark follows the "X" run).
        $bmRange = $cell.Range.Duplicate
        $d.Bookmarks.Add("_GoBack", $bmRange)

        $cell2 = $t.Cell($c.Row, $c.Col)
        $cell2.Range.Text = "X"
    } else {
        $cell.Range.Text = "X"
    }
}

Write-Output "Filled C/D transaction cells and relocated _GoBack bookmark."
